$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("level1")
$ws2 = $wb.Worksheets.Item("level2")
$ws3 = $wb.Worksheets.Item("level3")

# --- level1 (sheet1) ---
# vary_with_level for gw_depth row now TRUE
$ws1.Range("C16").Value2 = $true

# add new random variable row: theta_rake
$ws1.Range("A29").Value2 = "theta_rake"
$ws1.Range("B29").Value2 = $true
$ws1.Range("C29").Value2 = $false
$ws1.Range("D29").Value2 = "pipe-fault dip angle"
$ws1.Range("E29").Value2 = "deg"
$ws1.Range("F29").Value2 = "depends"
$ws1.Range("G29").Value2 = 10
$ws1.Range("I29").Value2 = -360
$ws1.Range("J29").Value2 = 360
$ws1.Range("K29").Value2 = "normal"

# --- level2 (sheet2) ---
$ws2.Range("C16").Value2 = $true

$ws2.Range("A29").Value2 = "theta_rake"
$ws2.Range("B29").Value2 = $true
$ws2.Range("C29").Value2 = $false
$ws2.Range("D29").Value2 = "pipe-fault dip angle"
$ws2.Range("E29").Value2 = "deg"
$ws2.Range("F29").Value2 = "depends"
$ws2.Range("G29").Value2 = 10
$ws2.Range("I29").Value2 = -360
$ws2.Range("J29").Value2 = 360
$ws2.Range("K29").Value2 = "normal"

# --- level3 (sheet3) ---
$ws3.Range("C16").Value2 = $true
$ws3.Range("F16").Value2 = "user provided"

$ws3.Range("A29").Value2 = "theta_rake"
$ws3.Range("B29").Value2 = $true
$ws3.Range("C29").Value2 = $false
$ws3.Range("D29").Value2 = "pipe-fault dip angle"
$ws3.Range("E29").Value2 = "deg"
$ws3.Range("F29").Value2 = "depends"
$ws3.Range("G29").Value2 = 10
$ws3.Range("I29").Value2 = -360
$ws3.Range("J29").Value2 = 360
$ws3.Range("K29").Value2 = "normal"

# --- selections / active sheet ---
$ws1.Range("G29").Select()

$ws2.Activate()
$ws2.Range("G29").Select()

$ws3.Activate()
$ws3.Range("F15,F16").Select()
$ws3.Range("F16").Activate()
